$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "session" column header
$ws.Range("D1").Value = "session"

# Morning session: rows 2-9 (time_start 8:00 through 11:30)
$ws.Range("D2:D9").Value = "Morning"

# Afternoon session: rows 10-17 (time_start 13:30 through 17:00)
$ws.Range("D10:D17").Value = "Afternoon"

# Move selection to the next empty row, as Excel does after data entry
$ws.Range("D18").Select()
